$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(203957296, Omri Ben Shabat: 7,8)"
$ws.Range("B1").Value = "(206532695, Matan Vakrat: -8,5)"
$ws.Range("C1").Value = "(302962915, Asher  Odeh: -8,-8)"
$ws.Range("D1").Value = "(308035542, Anastasia  Kubi: 6,-2)"
$ws.Range("E1").Value = "(311177802, Christina  Uksusman: -7,-7)"
$ws.Range("F1").Value = "(305251175, Or  Leder: 0,-10)"

$ws.Range("A3").Value = "cost: 812.7830972132518"
$ws.Range("A4").Value = "time: 98.47288715165648"
